$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.103.74"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.550.96"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "586.35"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").Value = "147.42"
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").Value = "27.53"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "3.005.16"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "63.029.05"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "2.552.93"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").Value = "337.02"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "65.80"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "1.62"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "8.39"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").Value = "  +8.88%  "
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +6.13%  "
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").Value = "178.62"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "420.01"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "0.401"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "39.76"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").Value = "150.39"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "20.89"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("E51").Value = "  -0.32%  "
